$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 997
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()

$ws.Range("H33").Value = 802.7
$ws.Range("I33").Value = 283.5
$ws.Range("K33").Value = 283.5
$ws.Range("M33").Value = -54.5

$ws.Range("H88").Value = 2893.2727
$ws.Range("I88").Value = 2487
$ws.Range("J88").Value = 3231.8333
$ws.Range("K88").Value = 2487
$ws.Range("L88").Value = 3231.8333
$ws.Range("M88").Value = -2081
$ws.Range("N88").Value = -4043.8333

$ws.Range("H91").Value = 2893.2727
$ws.Range("I91").Value = 2487
$ws.Range("J91").Value = 3231.8333
$ws.Range("K91").Value = 2487
$ws.Range("L91").Value = 3231.8333
$ws.Range("M91").Value = -1083
$ws.Range("N91").Value = -6039.8333

$ws.Range("H98").Value = 3155525.2
$ws.Range("I98").Value = 3498280
$ws.Range("J98").Value = 1670254.9
$ws.Range("K98").Value = 3498280
$ws.Range("L98").Value = 1670254.9
$ws.Range("M98").Value = -3496782
$ws.Range("N98").Value = -1673250.9

$ws.Range("H103").Value = 500
$ws.Range("J103").Value = 500
$ws.Range("L103").Value = 1500
$ws.Range("N103").Value = -2672

$ws.Range("H115").Value = 63607524
$ws.Range("I115").Value = 71557220
$ws.Range("K115").Value = 214671660
$ws.Range("M115").Value = -214670093

$ws.Range("H122").Value = 3155525.2
$ws.Range("I122").Value = 3498280
$ws.Range("J122").Value = 1670254.9
$ws.Range("K122").Value = 10494840
$ws.Range("L122").Value = 5010764.699999999
$ws.Range("M122").Value = -10492390
$ws.Range("N122").Value = -5015664.699999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 8000
$ws.Range("I8").Value = 8000
$ws.Range("K8").Value = 8000
$ws.Range("M8").Value = -7856

$ws.Range("H32").Value = 6389.6963
$ws.Range("I32").Value = 6389.6963
$ws.Range("K32").Value = 6389.6963
$ws.Range("M32").Value = -6102.6963

$ws.Range("H38").Value = 49000
$ws.Range("I38").Value = 49000
$ws.Range("K38").Value = 49000
$ws.Range("M38").Value = -48533

$ws.Range("H61").Value = 4515.6
$ws.Range("I61").Value = 3157
$ws.Range("J61").Value = 9950
$ws.Range("K61").Value = 3157
$ws.Range("L61").Value = 9950
$ws.Range("M61").Value = -2945
$ws.Range("N61").Value = -10374

$ws.Range("H74").Value = 125738.875
$ws.Range("I74").Value = 125738.875
$ws.Range("K74").Value = 125738.875
$ws.Range("M74").Value = -124864.875

$ws.Range("H77").Value = 125738.875
$ws.Range("I77").Value = 125738.875
$ws.Range("K77").Value = 628694.375
$ws.Range("M77").Value = -624326.375

$ws.Range("H122").Value = 3500
$ws.Range("I122").Value = 3500
$ws.Range("K122").Value = 10500
$ws.Range("M122").Value = -8050

$ws.Range("H136").Value = 4515.6
$ws.Range("I136").Value = 3157
$ws.Range("J136").Value = 9950
$ws.Range("K136").Value = 9471
$ws.Range("L136").Value = 29850
$ws.Range("M136").Value = -6921
$ws.Range("N136").Value = -34950

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1182.4375
$ws.Range("I20").Value = 1097.4783
$ws.Range("J20").Value = 1399.5555
$ws.Range("K20").Value = 1097.4783
$ws.Range("L20").Value = 1399.5555
$ws.Range("M20").Value = -850.4783
$ws.Range("N20").Value = -1893.5555

$ws.Range("H86").Value = 4300
$ws.Range("J86").Value = 5666.6665
$ws.Range("L86").Value = 5666.6665
$ws.Range("N86").Value = -7912.6665

$ws.Range("H89").Value = 4300
$ws.Range("J89").Value = 5666.6665
$ws.Range("L89").Value = 28333.3325
$ws.Range("N89").Value = -39565.3325

$ws.Range("H134").Value = 2432.9722
$ws.Range("I134").Value = 1926.9
$ws.Range("J134").Value = 4963.3335
$ws.Range("K134").Value = 5780.700000000001
$ws.Range("L134").Value = 14890.0005
$ws.Range("M134").Value = -3245.700000000001
$ws.Range("N134").Value = -19960.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3367337.2
$ws.Range("J99").Value = 100000
$ws.Range("L99").Value = 100000
$ws.Range("N99").Value = -102996

$ws.Range("H107").Value = 4270.763
$ws.Range("I107").Value = 574
$ws.Range("K107").Value = 574
$ws.Range("M107").Value = 1346

$ws.Range("H126").Value = 3367337.2
$ws.Range("J126").Value = 100000
$ws.Range("L126").Value = 300000
$ws.Range("N126").Value = -304940

$ws.Range("H132").Value = 2806.625
$ws.Range("I132").Value = 2925.5
$ws.Range("K132").Value = 8776.5
$ws.Range("M132").Value = -6246.5

$ws.Range("H134").Value = 3852.4187
$ws.Range("I134").Value = 4199.4326
$ws.Range("K134").Value = 12598.2978
$ws.Range("M134").Value = -10063.2978

$ws.Range("H141").Value = 128416.805
$ws.Range("J141").Value = 128416.805
$ws.Range("L141").Value = 128416.805
$ws.Range("N141").Value = -138776.805

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7681.4287
$ws.Range("I70").Value = 8439.691999999999
$ws.Range("K70").Value = 8439.691999999999
$ws.Range("M70").Value = -8169.691999999999

$ws.Range("H73").Value = 7681.4287
$ws.Range("I73").Value = 8439.691999999999
$ws.Range("K73").Value = 8439.691999999999
$ws.Range("M73").Value = -7503.691999999999

$ws.Range("H80").Value = 3521.889
$ws.Range("I80").Value = 3454.4688
$ws.Range("J80").Value = 4061.25
$ws.Range("K80").Value = 3454.4688
$ws.Range("L80").Value = 4061.25
$ws.Range("M80").Value = -2456.4688
$ws.Range("N80").Value = -6057.25

$ws.Range("H83").Value = 3521.889
$ws.Range("I83").Value = 3454.4688
$ws.Range("J83").Value = 4061.25
$ws.Range("K83").Value = 17272.344
$ws.Range("L83").Value = 20306.25
$ws.Range("M83").Value = -12280.344
$ws.Range("N83").Value = -30290.25

$ws.Range("H126").Value = 7876.778
$ws.Range("I126").Value = 15332.667
$ws.Range("J126").Value = 4148.8335
$ws.Range("K126").Value = 45998.001
$ws.Range("L126").Value = 12446.5005
$ws.Range("M126").Value = -43528.001
$ws.Range("N126").Value = -17386.5005

$ws.Range("H132").Value = 36449.152
$ws.Range("I132").Value = 40241.086
$ws.Range("J132").Value = 3269.75
$ws.Range("K132").Value = 120723.258
$ws.Range("L132").Value = 9809.25
$ws.Range("M132").Value = -118193.258
$ws.Range("N132").Value = -14869.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7075.8335
$ws.Range("I7").Value = 7110.5674
$ws.Range("J7").Value = 6818.8
$ws.Range("K7").Value = 7110.5674
$ws.Range("L7").Value = 6818.8
$ws.Range("M7").Value = -6998.5674
$ws.Range("N7").Value = -7042.8

$ws.Range("H61").Value = 1791
$ws.Range("I61").Value = 1814.4
$ws.Range("J61").Value = 1732.5
$ws.Range("K61").Value = 1814.4
$ws.Range("L61").Value = 1732.5
$ws.Range("M61").Value = -1612.4
$ws.Range("N61").Value = -2136.5

$ws.Range("H100").Value = 2000
$ws.Range("I100").Value = 2000
$ws.Range("K100").Value = 2000
$ws.Range("M100").Value = -1459

$ws.Range("H113").Value = 1791
$ws.Range("I113").Value = 1814.4
$ws.Range("J113").Value = 1732.5
$ws.Range("K113").Value = 1814.4
$ws.Range("L113").Value = 1732.5
$ws.Range("M113").Value = 355.5999999999999
$ws.Range("N113").Value = -6072.5

$ws.Range("H126").Value = 7075.8335
$ws.Range("I126").Value = 7110.5674
$ws.Range("J126").Value = 6818.8
$ws.Range("K126").Value = 21331.7022
$ws.Range("L126").Value = 20456.4
$ws.Range("M126").Value = -18861.7022
$ws.Range("N126").Value = -25396.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 172383
$ws.Range("I96").Value = 339766.34
$ws.Range("J96").Value = 4999.6665
$ws.Range("K96").Value = 339766.34
$ws.Range("L96").Value = 4999.6665
$ws.Range("M96").Value = -338393.34
$ws.Range("N96").Value = -7745.6665

$ws.Range("H132").Value = 2622.83
$ws.Range("I132").Value = 2916.0588
$ws.Range("K132").Value = 8748.1764
$ws.Range("M132").Value = -6218.1764

$ws.Range("H136").Value = 297457.25
$ws.Range("I136").Value = 348499.9
$ws.Range("J136").Value = 1410
$ws.Range("K136").Value = 1045499.7
$ws.Range("L136").Value = 4230
$ws.Range("M136").Value = -1042949.7
$ws.Range("N136").Value = -9330
